# Insert a new data row at row 197 (pushing the existing rows 197:327 down to
# 198:328) and populate the new row with this week's price-report figures for
# "Feria Lagunitas de Puerto Montt - Apio".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 197:327 down by one row, duplicating the formatting of the row
# above (matches Excel's default "insert copied cells" behaviour for a
# whole-row insert, which is what keeps the D column's date-number style).
$ws.Rows("197:197").Insert()

# Populate the newly-inserted row with the new observation.
$ws.Range("A197").Value = 4
$ws.Range("B197").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C197").Value = "Los Lagos"
$ws.Range("D197").Value = 44827
$ws.Range("E197").Value = 10
$ws.Range("F197").Value = 100112017
$ws.Range("G197").Value = "Apio"
$ws.Range("H197").Value = "Americana (o)"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 50
$ws.Range("K197").Value = 14500
$ws.Range("L197").Value = 15000
$ws.Range("M197").Value = 14750
$ws.Range("N197").Value = "`$/docena de matas"
$ws.Range("O197").Value = "Región de Coquimbo"
$ws.Range("P197").Value = 2458
$ws.Range("Q197").Value = 6
$ws.Range("R197").Value = "Hortaliza"
